# Update "想去人数" (want-to-go count) figures in column F across the
# four sheets of the 上海-漫展信息 workbook, reflecting a refreshed data
# scrape (gh-pages output regenerated at commit 456a3b4).

$wb = $excel.ActiveWorkbook

function Set-CellValue {
    param(
        [string]$SheetName,
        [string]$CellRef,
        [double]$NewValue
    )
    $ws = $wb.Worksheets.Item($SheetName)
    $ws.Range($CellRef).Value = $NewValue
}

# 展览 (Exhibitions) sheet
Set-CellValue "展览" "F7" 240
Set-CellValue "展览" "F10" 6938
Set-CellValue "展览" "F13" 356
Set-CellValue "展览" "F17" 2246
Set-CellValue "展览" "F18" 1499
Set-CellValue "展览" "F19" 654
Set-CellValue "展览" "F23" 182
Set-CellValue "展览" "F26" 1743
Set-CellValue "展览" "F36" 19
Set-CellValue "展览" "F38" 2734
Set-CellValue "展览" "F39" 75
Set-CellValue "展览" "F48" 11

# 演出 (Performances) sheet
Set-CellValue "演出" "F20" 56
Set-CellValue "演出" "F23" 479

# 本地生活 (Local life) sheet
Set-CellValue "本地生活" "F6" 1693
Set-CellValue "本地生活" "F8" 2739
Set-CellValue "本地生活" "F9" 1024
Set-CellValue "本地生活" "F10" 939
Set-CellValue "本地生活" "F14" 7374

# 全部类型 (All types) sheet
Set-CellValue "全部类型" "F6" 1693
Set-CellValue "全部类型" "F8" 2739
Set-CellValue "全部类型" "F9" 6938
Set-CellValue "全部类型" "F10" 1024
Set-CellValue "全部类型" "F12" 356
Set-CellValue "全部类型" "F16" 2246
Set-CellValue "全部类型" "F17" 1499
Set-CellValue "全部类型" "F22" 1743
Set-CellValue "全部类型" "F31" 56
Set-CellValue "全部类型" "F34" 479
Set-CellValue "全部类型" "F37" 19
Set-CellValue "全部类型" "F39" 2734
Set-CellValue "全部类型" "F40" 75
